$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '42.904.86'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '2.282.32'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue $ws.Range("D5") '249.95'
$ws.Range("E5").Value = '  -0.95%  '
Set-TextValue $ws.Range("D6") '0.644'
$ws.Range("E6").Value = '  +0.93%  '
Set-TextValue $ws.Range("D7") '79.19'
$ws.Range("E7").Value = '  +8.97%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.05%  '
Set-TextValue $ws.Range("D10") '40.98'
$ws.Range("E10").Value = '  +5.92%  '
Set-TextValue $ws.Range("D11") '0.0975'
$ws.Range("E11").Value = '  -0.30%  '
Set-TextValue $ws.Range("D12") '7.36'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").Value = '2.623.04'
$ws.Range("E14").Value = '  +0.25%  '
Set-TextValue $ws.Range("D15") '15.12'
$ws.Range("E15").Value = '  +0.68%  '
Set-TextValue $ws.Range("D16") '0.871'
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = '2.271.26'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '42.827.15'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").Value = '0.0₃0996'
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("E20").Value = '  -1.76%  '
Set-TextValue $ws.Range("D21") '72.19'
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("E23").Value = '  +1.14%  '
Set-TextValue $ws.Range("D24") '3.81'
$ws.Range("E24").Value = '  -2.31%  '
$ws.Range("E25").Value = '  -0.01%  '
Set-TextValue $ws.Range("D26") '11.34'
$ws.Range("E26").Value = '  -2.27%  '
$ws.Range("E27").Value = '  -4.31%  '
Set-TextValue $ws.Range("D28") '2.18'
$ws.Range("E28").Value = '  +2.53%  '
Set-TextValue $ws.Range("D29") '167.79'
$ws.Range("E29").Value = '  -0.37%  '
Set-TextValue $ws.Range("D30") '20.92'
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("E31").Value = '  +1.06%  '
Set-TextValue $ws.Range("D32") '0.0854'
$ws.Range("E32").Value = '  +6.22%  '
Set-TextValue $ws.Range("D33") '0.123'
$ws.Range("E33").Value = '  -5.58%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D34") '0.128'
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D35") '30.11'
$ws.Range("E35").Value = '  -2.19%  '
Set-TextValue $ws.Range("D36") '4.56'
$ws.Range("E36").Value = '  -2.60%  '
Set-TextValue $ws.Range("D37") '4.81'
$ws.Range("E37").Value = '  +0.70%  '
Set-TextValue $ws.Range("D38") '0.0304'
$ws.Range("E38").Value = '  -2.40%  '
Set-TextValue $ws.Range("D39") '13.70'
$ws.Range("E39").Value = '  +3.30%  '
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D42") '0.210'
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D43") '112.08'
$ws.Range("E43").Value = '  +15.88%  '
Set-TextValue $ws.Range("D44") '61.29'
$ws.Range("E44").Value = '  -0.65%  '
Set-TextValue $ws.Range("D45") '8.92'
$ws.Range("E45").Value = '  -3.00%  '
Set-TextValue $ws.Range("D47") '4.64'
$ws.Range("E47").Value = '  -6.87%  '
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("E49").Value = '  -2.92%  '
Set-TextValue $ws.Range("D50") '1.17'
$ws.Range("E50").Value = '  -2.52%  '
Set-TextValue $ws.Range("D51") '4.25'
$ws.Range("E51").Value = '  +0.01%  '
